$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")

# --- Hoja2: extend the validation source list with 3 new entries ---
$ws2.Range("A4").Value = "Edecán"
$ws2.Range("A5").Value = "Limpieza"
$ws2.Range("A6").Value = "Otros"

# --- Hoja1: widen the data validation on column D to use the new list range ---
$ws1.Range("D2:D1300").Validation.Delete()
$ws1.Range("D2:D1300").Validation.Add(3, 1, 1, "=Hoja2!`$A`$2:`$A`$6")

# --- Hoja1: page setup (paper size / orientation) ---
$ps = $ws1.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

# --- Hoja1: scroll the view down near the bottom of the used data and select H1307 ---
$ws1.Activate()
$ws1.Range("H1307").Select()
$excel.ActiveWindow.ScrollRow = 1299
$excel.ActiveWindow.ScrollColumn = 1

# --- Hoja2: select C8 ---
$ws2.Activate()
$ws2.Range("C8").Select()

# Re-activate Hoja1 as the selected/visible sheet (matches tabSelected="1" on sheet1)
$ws1.Activate()
